$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column H (before the existing "url" column), shifting url to I
$ws.Columns("H").Insert()

# Add the new row of data (arrowtooth flounder) in row 7 -- entered in the
# same order the original author typed it: filename, type, species, common
# name, then the url, and finally the new column header.
$ws.Range("A7").Value = 10110
$ws.Range("D7").Value = "Data_ATF2021_1993plus_all data.csv"
$ws.Range("E7").Value = "CSV"
$ws.Range("B7").Value = "Atheresthes stomias"
$ws.Range("C7").Value = "arrowtooth flounder"
$ws.Range("F7").Value = 1993
$ws.Range("G7").Value = 2021
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = "https://drive.google.com/drive/folders/1oC3-kAIcyy7-W2joSkJwX66ynXfo_hCc"

# Turn the URL in I7 into a real hyperlink (keeps display text = URL)
$ws.Hyperlinks.Add($ws.Range("I7"), "https://drive.google.com/drive/folders/1oC3-kAIcyy7-W2joSkJwX66ynXfo_hCc")

# Header + values for the new "area_units_correct" column
$ws.Range("H1").Value = "area_units_correct"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0.001

# Column widths to match the new layout
$ws.Columns("C").ColumnWidth = 18.1796875
$ws.Columns("D").ColumnWidth = 35.81640625

# Restore view to top-left and select column A:B
$ws.Range("A1").Select()
$ws.Range("A1:B1048576").Select()

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1
